$d = $word.ActiveDocument

# The document used to end with a short "Requisitos" section followed by
# three site-chrome paragraphs that came from the Jekyll page template:
#   (blank paragraph)
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#    pages. Original theme under Creative Commons Attribution"
# On this site rebuild that chrome is gone. Locate the "Ver no Jupiter ..."
# paragraph, then delete it together with the blank paragraph right
# before it and the copyright paragraph right after it, leaving the
# "LOQ4031: ..." requirement paragraph (and everything that follows the
# copyright line) untouched.

$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $searchRange.Find.Found) {
    $jupPara = $searchRange.Paragraphs(1)
    $prevPara = $jupPara.Previous()
    $nextPara = $jupPara.Next()

    $delStart = $prevPara.Range.Start
    $delEnd = $nextPara.Range.End

    $d.Range($delStart, $delEnd).Delete()
}
